$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "19.09.2018, Wed"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = "`n"

# Writing a literal newline makes the host measure/auto-size the row; the
# source workbook keeps this row at the sheet's default (non-custom) height,
# same as the other date rows (13/15/18) that already hold a "`n" in column K.
$ws.Rows(20).AutoFit()
